# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1207
$ws1.Range("F5").Value = 338
$ws1.Range("F6").Value = 308
$ws1.Range("F7").Value = 3414
$ws1.Range("F8").Value = 233
$ws1.Range("F9").Value = 730
$ws1.Range("F10").Value = 908
$ws1.Range("F11").Value = 323
$ws1.Range("F16").Value = 1998
$ws1.Range("F19").Value = 8
$ws1.Range("F22").Value = 264

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 292

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 292
$ws4.Range("F13").Value = 1207
$ws4.Range("F14").Value = 338
$ws4.Range("F18").Value = 308
$ws4.Range("F19").Value = 3415
$ws4.Range("F21").Value = 233
$ws4.Range("F25").Value = 730
$ws4.Range("F26").Value = 908
$ws4.Range("F27").Value = 323
$ws4.Range("F34").Value = 1998
$ws4.Range("F39").Value = 8
$ws4.Range("F49").Value = 264

$wb.Save()
